{"js": "// Clean up the duplicated \"egXML\" / \"egXMLTable\" gloss-list paragraph\n// styles that had accumulated in the stylesheet (egXML, egXML0..egXML4,\n// egXMLTable, egXMLTable0..egXMLTable4). Keep exactly one definition of\n// each (re-created with identical formatting), and drop the rest.\n\nconst styles = context.document.getStyles();\nstyles.load(\"items/nameLocal,items/type\");\nawait context.sync();\n\n// Collect the indices (document order) of every custom paragraph style\n// named \"egXML\" or \"egXMLTable\" - these are the duplicated gloss-example\n// styles being cleaned up.\nconst dupIndices = [];\nfor (let i = 0; i < styles.items.length; i++) {\n  const nm = styles.items[i].nameLocal;\n  const ty = styles.items[i].type;\n  if ((nm === \"egXML\" || nm === \"egXMLTable\") && ty === Word.StyleType.paragraph) {\n    dupIndices.push(i);\n  }\n}\n\n// Delete them all (highest index first so earlier indices stay valid as\n// the collection shrinks).\ndupIndices.sort((a, b) => b - a);\nfor (const idx of dupIndices) {\n  styles.items[idx].delete();\n}\nawait context.sync();\n\n// Re-create a single canonical copy of each style, with the same\n// formatting every duplicate shared: based on Normal, Courier font,\n// quick-style, egXML at 10pt and egXMLTable at 9pt with 4pt space-before.\nconst egXML = context.document.addStyle(\"egXML\", Word.StyleType.paragraph);\nawait context.sync();\n\nconst egXMLStyle = context.document.getStyles().getByName(\"egXML\");\negXMLStyle.baseStyle = \"Normal\";\negXMLStyle.quickStyle = true;\negXMLStyle.font.name = \"Courier\";\negXMLStyle.font.size = 10;\nawait context.sync();\n\nconst egXMLTable = context.document.addStyle(\"egXMLTable\", Word.StyleType.paragraph);\nawait context.sync();\n\nconst egXMLTableStyle = context.document.getStyles().getByName(\"egXMLTable\");\negXMLTableStyle.baseStyle = \"Normal\";\negXMLTableStyle.quickStyle = true;\negXMLTableStyle.font.name = \"Courier\";\negXMLTableStyle.font.size = 9;\negXMLTableStyle.paragraphFormat.spaceBefore = 4;\nawait context.sync();\n", "ps1": "# Clean up the duplicated \"egXML\" / \"egXMLTable\" gloss-list paragraph\n# styles that had accumulated in the stylesheet (egXML, egXML0..egXML4,\n# egXMLTable, egXMLTable0..egXMLTable4). Keep exactly one definition of\n# each (re-created with identical formatting), and drop the rest.\n\n$d = $word.ActiveDocument\n\n# Find every custom paragraph style named \"egXML\" or \"egXMLTable\" - these\n# are the duplicated gloss-example styles being cleaned up.\n# wdStyleTypeParagraph = 1\n$indices = @()\nfor ($i = 1; $i -le $d.Styles.Count; $i++) {\n    $s = $d.Styles.Item($i)\n    if (($s.NameLocal -eq \"egXML\" -or $s.NameLocal -eq \"egXMLTable\") -and $s.Type -eq 1) {\n        $indices += $i\n    }\n}\n\n# Delete them all (highest index first so earlier indices stay valid as\n# the collection shrinks).\n$sorted = $indices | Sort-Object -Descending\nforeach ($i in $sorted) {\n    $d.Styles.Item($i).Delete()\n}\n\n# Re-create a single canonical copy of each style, with the same\n# formatting every duplicate shared: based on Normal, Courier font,\n# quick-style, egXML at 10pt and egXMLTable at 9pt with 4pt space-before.\n$egXML = $d.Styles.Add(\"egXML\", 1)\n$egXML.BaseStyle = $d.Styles.Item(\"Normal\")\n$egXML.QuickStyle = $true\n$egXML.Font.Name = \"Courier\"\n$egXML.Font.Size = 10\n\n$egXMLTable = $d.Styles.Add(\"egXMLTable\", 1)\n$egXMLTable.BaseStyle = $d.Styles.Item(\"Normal\")\n$egXMLTable.QuickStyle = $true\n$egXMLTable.Font.Name = \"Courier\"\n$egXMLTable.Font.Size = 9\n$egXMLTable.ParagraphFormat.SpaceBefore = 4\n"}
